# Replace the sample "login credentials" sheet with Create-SKU validation
# data used by the automation test resource.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the old mailto: hyperlink on C2, then wipe the old sample table.
$ws.Range("C2").Hyperlinks.Delete()
$ws.Cells.Clear()

# Columns B holding "05"/"120"/"2" must stay text (leading zeros / strings),
# not get auto-coerced to numbers, so pre-format them before writing.
$ws.Range("B7:B8").NumberFormat = "@"
$ws.Range("B12:B13").NumberFormat = "@"

$data = @(
    @("HSN Key", 1245),
    @("SKU Name", "test SKU"),
    @("SKU Local Name", "test SKU Local"),
    @("Description", "This is test SKU"),
    @("Category", "Biscuits And Choclates"),
    @("sub-Category", "Biscuits"),
    @("GST Rate", "05"),
    @("CESS Rate", "05"),
    @("Amount Type", "Net"),
    @("Brand Name", "Parle Agro"),
    @("variation Name", "test SKU automation"),
    @("Price", "120"),
    @("value", "2")
)

$r = 1
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# Clear the temporary "@" number formats back off so no stray style index is
# left on the cells (Normal == default style, xfId 0).
$ws.UsedRange.Style = "Normal"

$ws.Columns.Item(1).ColumnWidth = 21.42578125
$ws.Columns.Item(2).ColumnWidth = 19.85546875
$ws.Columns.Item(3).ColumnWidth = 15
$ws.Columns.Item(4).ColumnWidth = 15.42578125

$ws.Range("A1:C15").Select()
$ws.Application.ActiveWindow.ScrollRow = 7
